$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 16:44"

# Swap Croacia/Birmania ranking (row 89/90) - country names
$ws.Cells.Item(89, 1).Value = "Birmania"
$ws.Cells.Item(90, 1).Value = "Croacia"

# Swap Islas Malvinas/Montserrat ranking (row 215/216) - country names
$ws.Cells.Item(215, 1).Value = "Montserrat"
$ws.Cells.Item(216, 1).Value = "Islas Malvinas"

# Updated statistics snapshot (new scrape values)
# Row 4
$ws.Cells.Item(4, 2).Value = 7641804
$ws.Cells.Item(4, 3).Value = 4892
$ws.Cells.Item(4, 4).Value = 4851087
$ws.Cells.Item(4, 5).Value = 2576072
$ws.Cells.Item(4, 7).Value = 34
$ws.Cells.Item(4, 8).Value = 214645

# Row 5
$ws.Cells.Item(5, 2).Value = 6644087
$ws.Cells.Item(5, 3).Value = 21907
$ws.Cells.Item(5, 4).Value = 5611050
$ws.Cells.Item(5, 5).Value = 930068
$ws.Cells.Item(5, 7).Value = 255
$ws.Cells.Item(5, 8).Value = 102969

# Row 26
$ws.Cells.Item(26, 2).Value = 302509
$ws.Cells.Item(26, 3).Value = 938
$ws.Cells.Item(26, 5).Value = 31003
$ws.Cells.Item(26, 7).Value = 4
$ws.Cells.Item(26, 8).Value = 9606

# Row 70
$ws.Cells.Item(70, 2).Value = 41957
$ws.Cells.Item(70, 3).Value = 459
$ws.Cells.Item(70, 4).Value = 35182
$ws.Cells.Item(70, 5).Value = 6436
$ws.Cells.Item(70, 7).Value = 9
$ws.Cells.Item(70, 8).Value = 339

# Row 72
$ws.Cells.Item(72, 2).Value = 39449
$ws.Cells.Item(72, 3).Value = 22
$ws.Cells.Item(72, 5).Value = 13059

# Row 89
$ws.Cells.Item(89, 2).Value = 18781
$ws.Cells.Item(89, 3).Value = 987
$ws.Cells.Item(89, 4).Value = 5548
$ws.Cells.Item(89, 5).Value = 12789
$ws.Cells.Item(89, 7).Value = 32
$ws.Cells.Item(89, 8).Value = 444

# Row 90
$ws.Cells.Item(90, 2).Value = 17797
$ws.Cells.Item(90, 3).Value = 138
$ws.Cells.Item(90, 4).Value = 16031
$ws.Cells.Item(90, 5).Value = 1466
$ws.Cells.Item(90, 7).Value = 2
$ws.Cells.Item(90, 8).Value = 300

# Row 94
$ws.Cells.Item(94, 2).Value = 15089
$ws.Cells.Item(94, 3).Value = 37
$ws.Cells.Item(94, 4).Value = 14306
$ws.Cells.Item(94, 5).Value = 449
$ws.Cells.Item(94, 7).Value = 1
$ws.Cells.Item(94, 8).Value = 334

# Row 95
$ws.Cells.Item(95, 2).Value = 14527
$ws.Cells.Item(95, 3).Value = 70
$ws.Cells.Item(95, 5).Value = 3062

# Row 96
$ws.Cells.Item(96, 2).Value = 14410
$ws.Cells.Item(96, 3).Value = 144
$ws.Cells.Item(96, 4).Value = 8825
$ws.Cells.Item(96, 5).Value = 5185
$ws.Cells.Item(96, 7).Value = 4
$ws.Cells.Item(96, 8).Value = 400

# Row 107
$ws.Cells.Item(107, 2).Value = 9974
$ws.Cells.Item(107, 3).Value = 39
$ws.Cells.Item(107, 4).Value = 8794
$ws.Cells.Item(107, 5).Value = 1102

# Row 143
$ws.Cells.Item(143, 2).Value = 3483
$ws.Cells.Item(143, 3).Value = 81
$ws.Cells.Item(143, 5).Value = 211

# Row 215
$ws.Cells.Item(215, 4).Value = 12
$ws.Cells.Item(215, 8).Value = 1

# Row 216
$ws.Cells.Item(216, 4).Value = 13
$ws.Cells.Item(216, 8).Value = 0

